$wb = $excel.ActiveWorkbook

# --- Rename header cells on existing sheets ---
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("B1").Value = "Weekly_PO_Qty"

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after "Monthly Trend" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "PO Forecast"

# --- Header row: copy formatting (bold, border, center/top align) from an
#     existing header cell, then overwrite with the new header text ---
$ws1.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- Column A (dates) formatting: copy the date-number-format style used
#     for column A of the existing sheets down the forecast rows ---
$ws1.Range("A2").Copy()
$ws3.Range("A2:A81").PasteSpecial(-4122)

# --- Bulk-write the forecast data (rows 2-81) ---
$data = New-Object 'object[,]' 80,4
$data[0,0] = 44934.99999999999; $data[0,1] = 121; $data[0,2] = -379.7055199950197; $data[0,3] = 655.257621091397
$data[1,0] = 44941.99999999999; $data[1,1] = 127; $data[1,2] = -350.4221492416565; $data[1,3] = 641.8243460915357
$data[2,0] = 44955.99999999999; $data[2,1] = 140; $data[2,2] = -363.7831841310563; $data[2,3] = 658.8698417500337
$data[3,0] = 44962.99999999999; $data[3,1] = 146; $data[3,2] = -347.2662132422518; $data[3,3] = 728.7554612087218
$data[4,0] = 44969.99999999999; $data[4,1] = 153; $data[4,2] = -382.8052668780924; $data[4,3] = 691.3620597647648
$data[5,0] = 44976.99999999999; $data[5,1] = 159; $data[5,2] = -371.4235908806613; $data[5,3] = 711.6239076986212
$data[6,0] = 44983.99999999999; $data[6,1] = 165; $data[6,2] = -345.4226383145283; $data[6,3] = 694.0570257261144
$data[7,0] = 44990.99999999999; $data[7,1] = 172; $data[7,2] = -360.7967397878007; $data[7,3] = 676.5394010057023
$data[8,0] = 44997.99999999999; $data[8,1] = 178; $data[8,2] = -373.679024498676; $data[8,3] = 663.213907322383
$data[9,0] = 45004.99999999999; $data[9,1] = 185; $data[9,2] = -331.5733541875674; $data[9,3] = 695.6911817972725
$data[10,0] = 45011.99999999999; $data[10,1] = 191; $data[10,2] = -346.8135640706672; $data[10,3] = 716.4652770562556
$data[11,0] = 45025.99999999999; $data[11,1] = 204; $data[11,2] = -329.1762719456173; $data[11,3] = 718.9431128977265
$data[12,0] = 45032.99999999999; $data[12,1] = 210; $data[12,2] = -344.2294330874971; $data[12,3] = 749.7098841330242
$data[13,0] = 45039.99999999999; $data[13,1] = 217; $data[13,2] = -274.4744604484112; $data[13,3] = 730.5135090851203
$data[14,0] = 45060.99999999999; $data[14,1] = 236; $data[14,2] = -272.8246601862285; $data[14,3] = 722.033479135216
$data[15,0] = 45067.99999999999; $data[15,1] = 242; $data[15,2] = -258.5841810463382; $data[15,3] = 775.7413649879302
$data[16,0] = 45074.99999999999; $data[16,1] = 249; $data[16,2] = -240.7544834355683; $data[16,3] = 729.2892916073264
$data[17,0] = 45081.99999999999; $data[17,1] = 255; $data[17,2] = -282.3074994523296; $data[17,3] = 742.7344719763053
$data[18,0] = 45088.99999999999; $data[18,1] = 261; $data[18,2] = -233.94736504047; $data[18,3] = 800.4580740889675
$data[19,0] = 45102.99999999999; $data[19,1] = 274; $data[19,2] = -246.2861979528329; $data[19,3] = 776.7892727111168
$data[20,0] = 45109.99999999999; $data[20,1] = 280; $data[20,2] = -238.1443330900238; $data[20,3] = 814.2893150439952
$data[21,0] = 45116.99999999999; $data[21,1] = 287; $data[21,2] = -232.284268190267; $data[21,3] = 765.8223398160549
$data[22,0] = 45123.99999999999; $data[22,1] = 293; $data[22,2] = -207.0975498644625; $data[22,3] = 799.4463746932246
$data[23,0] = 45130.99999999999; $data[23,1] = 300; $data[23,2] = -194.2693053249469; $data[23,3] = 832.1914639882709
$data[24,0] = 45137.99999999999; $data[24,1] = 306; $data[24,2] = -190.7774297808317; $data[24,3] = 786.1540051341814
$data[25,0] = 45144.99999999999; $data[25,1] = 312; $data[25,2] = -210.4085848795311; $data[25,3] = 838.4497669123662
$data[26,0] = 45151.99999999999; $data[26,1] = 319; $data[26,2] = -180.657719450781; $data[26,3] = 864.7556665783927
$data[27,0] = 45158.99999999999; $data[27,1] = 325; $data[27,2] = -210.4900795053829; $data[27,3] = 838.2359662914439
$data[28,0] = 45165.99999999999; $data[28,1] = 332; $data[28,2] = -155.4645022219708; $data[28,3] = 869.9472181343388
$data[29,0] = 45172.99999999999; $data[29,1] = 338; $data[29,2] = -145.2997811751423; $data[29,3] = 851.2513288811197
$data[30,0] = 45179.99999999999; $data[30,1] = 344; $data[30,2] = -195.8153416234419; $data[30,3] = 831.914046515911
$data[31,0] = 45186.99999999999; $data[31,1] = 351; $data[31,2] = -168.9922668008755; $data[31,3] = 856.7531304305359
$data[32,0] = 45193.99999999999; $data[32,1] = 357; $data[32,2] = -130.3056093760337; $data[32,3] = 896.5581617247296
$data[33,0] = 45200.99999999999; $data[33,1] = 364; $data[33,2] = -144.4410251047299; $data[33,3] = 934.1033449516706
$data[34,0] = 45207.99999999999; $data[34,1] = 370; $data[34,2] = -144.804089064971; $data[34,3] = 898.6933437798069
$data[35,0] = 45214.99999999999; $data[35,1] = 376; $data[35,2] = -123.7627974712771; $data[35,3] = 883.9414347705014
$data[36,0] = 45221.99999999999; $data[36,1] = 383; $data[36,2] = -140.2544609410534; $data[36,3] = 856.0460348759937
$data[37,0] = 45228.99999999999; $data[37,1] = 389; $data[37,2] = -140.2400036697561; $data[37,3] = 888.4667860204147
$data[38,0] = 45235.99999999999; $data[38,1] = 396; $data[38,2] = -69.44118728930165; $data[38,3] = 903.3722061048553
$data[39,0] = 45242.99999999999; $data[39,1] = 402; $data[39,2] = -132.2098017513901; $data[39,3] = 909.9618582120121
$data[40,0] = 45249.99999999999; $data[40,1] = 408; $data[40,2] = -144.5723298392905; $data[40,3] = 906.1253112547979
$data[41,0] = 45256.99999999999; $data[41,1] = 415; $data[41,2] = -135.1009403275689; $data[41,3] = 948.6816425673671
$data[42,0] = 45270.99999999999; $data[42,1] = 428; $data[42,2] = -97.28659588511357; $data[42,3] = 955.335156383103
$data[43,0] = 45298.99999999999; $data[43,1] = 453; $data[43,2] = -38.20799508065055; $data[43,3] = 969.7029921416888
$data[44,0] = 45312.99999999999; $data[44,1] = 466; $data[44,2] = -35.2621145073201; $data[44,3] = 987.4563292751451
$data[45,0] = 45319.99999999999; $data[45,1] = 472; $data[45,2] = -55.09969291784679; $data[45,3] = 987.9419820074181
$data[46,0] = 45326.99999999999; $data[46,1] = 479; $data[46,2] = -53.81041793355724; $data[46,3] = 986.7833738351499
$data[47,0] = 45333.99999999999; $data[47,1] = 485; $data[47,2] = -36.71616726810342; $data[47,3] = 987.2063335221197
$data[48,0] = 45347.99999999999; $data[48,1] = 498; $data[48,2] = -37.01383507949769; $data[48,3] = 1023.202651873367
$data[49,0] = 45361.99999999999; $data[49,1] = 511; $data[49,2] = 26.71813835465261; $data[49,3] = 1029.879749885634
$data[50,0] = 45368.99999999999; $data[50,1] = 517; $data[50,2] = 10.80100765458229; $data[50,3] = 1033.824286604836
$data[51,0] = 45375.99999999999; $data[51,1] = 523; $data[51,2] = 34.45442787151968; $data[51,3] = 1025.557269548818
$data[52,0] = 45382.99999999999; $data[52,1] = 530; $data[52,2] = 16.19743140786628; $data[52,3] = 1029.297571417621
$data[53,0] = 45389.99999999999; $data[53,1] = 536; $data[53,2] = 28.23622179852481; $data[53,3] = 1038.534442042763
$data[54,0] = 45403.99999999999; $data[54,1] = 549; $data[54,2] = 58.59782074659157; $data[54,3] = 1016.66450731664
$data[55,0] = 45410.99999999999; $data[55,1] = 555; $data[55,2] = 22.41080798420475; $data[55,3] = 1081.225454149257
$data[56,0] = 45417.99999999999; $data[56,1] = 562; $data[56,2] = 33.30443915644934; $data[56,3] = 1086.751868579961
$data[57,0] = 45424.99999999999; $data[57,1] = 568; $data[57,2] = 13.89471159189415; $data[57,3] = 1052.973231384031
$data[58,0] = 45431.99999999999; $data[58,1] = 575; $data[58,2] = 48.63678368380817; $data[58,3] = 1072.778198893454
$data[59,0] = 45438.99999999999; $data[59,1] = 581; $data[59,2] = 81.62657226949908; $data[59,3] = 1069.413784206111
$data[60,0] = 45445.99999999999; $data[60,1] = 587; $data[60,2] = 96.76314321180988; $data[60,3] = 1095.182518779825
$data[61,0] = 45459.99999999999; $data[61,1] = 600; $data[61,2] = 89.17829447901036; $data[61,3] = 1141.603278718213
$data[62,0] = 45473.99999999999; $data[62,1] = 613; $data[62,2] = 80.78858461021284; $data[62,3] = 1128.96484368241
$data[63,0] = 45529.99999999999; $data[63,1] = 664; $data[63,2] = 159.2804217779226; $data[63,3] = 1172.881452116466
$data[64,0] = 45543.99999999999; $data[64,1] = 677; $data[64,2] = 139.8559020726672; $data[64,3] = 1142.650116780687
$data[65,0] = 45550.99999999999; $data[65,1] = 683; $data[65,2] = 164.2680045350392; $data[65,3] = 1216.64655617228
$data[66,0] = 45557.99999999999; $data[66,1] = 690; $data[66,2] = 184.8515275649114; $data[66,3] = 1185.558549244272
$data[67,0] = 45564.99999999999; $data[67,1] = 696; $data[67,2] = 157.2310330980428; $data[67,3] = 1213.374961192546
$data[68,0] = 45571.99999999999; $data[68,1] = 702; $data[68,2] = 179.132280109335; $data[68,3] = 1228.239638907166
$data[69,0] = 45578.99999999999; $data[69,1] = 709; $data[69,2] = 196.1813022936483; $data[69,3] = 1168.996309873338
$data[70,0] = 45585.99999999999; $data[70,1] = 715; $data[70,2] = 232.0863609228464; $data[70,3] = 1241.968293396091
$data[71,0] = 45599.99999999999; $data[71,1] = 728; $data[71,2] = 202.689658358079; $data[71,3] = 1228.04651094101
$data[72,0] = 45606.99999999999; $data[72,1] = 734; $data[72,2] = 205.3499533245499; $data[72,3] = 1263.072909181906
$data[73,0] = 45613.99999999999; $data[73,1] = 741; $data[73,2] = 236.2288148015432; $data[73,3] = 1279.859908913702
$data[74,0] = 45620.99999999999; $data[74,1] = 747; $data[74,2] = 254.4619341065068; $data[74,3] = 1255.232213385711
$data[75,0] = 45627.99999999999; $data[75,1] = 754; $data[75,2] = 232.0059416344455; $data[75,3] = 1240.499953359372
$data[76,0] = 45634.99999999999; $data[76,1] = 760; $data[76,2] = 256.9583662765845; $data[76,3] = 1282.443008226728
$data[77,0] = 45641.99999999999; $data[77,1] = 766; $data[77,2] = 246.3255962396683; $data[77,3] = 1322.53040366085
$data[78,0] = 45648.99999999999; $data[78,1] = 773; $data[78,2] = 259.2508786574508; $data[78,3] = 1242.117040570454
$data[79,0] = 45655.99999999999; $data[79,1] = 779; $data[79,2] = 255.285463874972; $data[79,3] = 1288.033950510305
$ws3.Range("A2:D81").Value = $data

# --- Restore the originally-active sheet/tab selection ---
$ws1.Activate()
